# Update marksheet scores: correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct answers count 3 -> 5
$ws.Range("B11").Value = 5

# Total row: total marks 54 -> 90
$ws.Range("B12").Value = 90

# Total row: corr/total text "47/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
